$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 827, shifting existing rows 827-868 down to 828-869
$ws.Rows.Item(827).Insert()

# Populate the newly inserted row 827 with the new data.
# Force column A to stay as plain text (not get auto-converted to a date
# serial number) by setting the number format to Text before assigning the
# value, then clear the formatting again so the cell ends up with the same
# default (unstyled) look as the other date-text cells in the column.
$cellA = $ws.Cells.Item(827, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2026/02/16"
$cellA.ClearFormats()

$ws.Cells.Item(827, 2).Value = "月"
$ws.Cells.Item(827, 3).Value = 0
$ws.Cells.Item(827, 4).Value = 201
